# Daily attendance processing - 2025-11-23 12:39:09
# Reorders the comma-separated "Recorded By" names in column G so that
# duplicate/alias audit entries list the canonical "system" actor first
# (rotating the trailing/alias entry to the front of the list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of exact current text -> corrected text for the "Recorded By" column (G).
$map = @{
    "backup@backdoor.com, System, system" = "system, backup@backdoor.com, System";
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com";
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com";
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = $ws.UsedRange.Rows.Count }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Text

    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
